$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "52.166.26"
$ws.Range("E2").Value = "  +0.93%  "

# Row 3
$ws.Range("D3").Value = "2.877.03"
$ws.Range("E3").Value = "  +3.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'352.74"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6
$ws.Range("D6").Value = "'112.71"
$ws.Range("E6").Value = "  +3.55%  "

# Row 7
$ws.Range("E7").Value = "  +2.28%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("E9").Value = "  +2.87%  "

# Row 10
$ws.Range("D10").Value = "'40.61"
$ws.Range("E10").Value = "  +2.01%  "

# Row 11
$ws.Range("E11").Value = "  -0.63%  "

# Row 12
$ws.Range("D12").Value = "'0.0855"
$ws.Range("E12").Value = "  +2.12%  "

# Row 13
$ws.Range("D13").Value = "'20.25"
$ws.Range("E13").Value = "  +0.60%  "

# Row 14
$ws.Range("E14").Value = "  +2.26%  "

# Row 15
$ws.Range("D15").Value = "3.331.62"
$ws.Range("E15").Value = "  +3.42%  "

# Row 16
$ws.Range("D16").Value = "2.905.60"
$ws.Range("E16").Value = "  +3.97%  "

# Row 17
$ws.Range("D17").Value = "'0.992"
$ws.Range("E17").Value = "  +6.79%  "

# Row 18
$ws.Range("D18").Value = "52.163.77"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("D19").Value = "'3.41"
$ws.Range("E19").Value = "  +8.62%  "

# Row 20
$ws.Range("D20").Value = "'7.70"
$ws.Range("E20").Value = "  -0.44%  "

# Row 21
$ws.Range("D21").Value = "'13.71"
$ws.Range("E21").Value = "  +4.10%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0977"
$ws.Range("E22").Value = "  +1.27%  "

# Row 23
$ws.Range("D23").Value = "'70.81"
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").Value = "'270.94"
$ws.Range("E24").Value = "  +1.40%  "

# Row 25
$ws.Range("D25").Value = "'2.79"
$ws.Range("E25").Value = "  +2.16%  "

# Row 26
$ws.Range("D26").Value = "'26.59"
$ws.Range("E26").Value = "  +1.81%  "

# Row 27
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("E28").Value = "  +0.99%  "

# Row 29
$ws.Range("D29").Value = "'10.58"
$ws.Range("E29").Value = "  +3.33%  "

# Row 30
$ws.Range("D30").Value = "'38.91"
$ws.Range("E30").Value = "  +4.70%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'2.27"
$ws.Range("E31").Value = "  +1.98%  "

# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.31"
$ws.Range("E32").Value = "  +1.93%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "'52.67"
$ws.Range("E33").Value = "  +1.91%  "

# Row 34
$ws.Range("D34").Value = "'0.0455"
$ws.Range("E34").Value = "  +0.58%  "

# Row 35
$ws.Range("E35").Value = "  +8.30%  "

# Row 36
$ws.Range("E36").Value = "  +0.16%  "

# Row 37
$ws.Range("E37").Value = "  +0.04%  "

# Row 38
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "  +5.66%  "

# Row 39
$ws.Range("D39").Value = "'19.00"
$ws.Range("E39").Value = "  +2.79%  "

# Row 40
$ws.Range("E40").Value = "  +3.46%  "

# Row 41
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +2.05%  "

# Row 42
$ws.Range("E42").Value = "  +1.97%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'22.62"
$ws.Range("E43").Value = "  +2.37%  "

# Row 44
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'122.59"
$ws.Range("E44").Value = "  +1.90%  "

# Row 45
$ws.Range("E45").Value = "  +2.01%  "

# Row 46
$ws.Range("D46").Value = "'3.56"
$ws.Range("E46").Value = "  +7.68%  "

# Row 47
$ws.Range("D47").Value = "2.182.29"
$ws.Range("E47").Value = "  +2.71%  "

# Row 48
$ws.Range("E48").Value = "  +6.99%  "

# Row 49
$ws.Range("D49").Value = "'0.245"
$ws.Range("E49").Value = "  +16.24%  "

# Row 50
$ws.Range("D50").Value = "'0.964"
$ws.Range("E50").Value = "  +6.28%  "

# Row 51
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "'0.0322"
$ws.Range("E51").Value = "  +13.32%  "
